$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (A7): clarify which ftl file to use
$ws.Range("A7").Value = "Make a new tempate and paste the content of the downloaded ftl (functionality_tabs_v2.ftl) file inside the template"

# Row 4 (A4): remove the "Upload the downloaded ftl files to MOLGENIS" instruction entirely
$ws.Range("A4").ClearContents()

# Row 11 (A11): add ":test" after the last two values in the example input
$ws.Range("A11").Value = "Typ in this in the input box: chromosome6_a_c:test,chromome6_d_h:test,chromome_i_L:test,chromosome6_array:test"

# Update the sheet view: no frozen/scrolled top-left cell, and selection now on A7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A7").Select()
